# Updated symbol list on Sat Jan 21 06:55:35 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for the cryptos sheet
# with the latest scraped quotes. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the sheet's
# existing inline-string / General-format cells) instead of auto-coercing
# them to numbers or percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.12"
$ws.Range("E2").Value = "'5.51%"
$ws.Range("D3").Value = "'34.87"
$ws.Range("E3").Value = "'12.57%"
$ws.Range("D4").Value = "'5.129"
$ws.Range("E4").Value = "'4.16%"
$ws.Range("D5").Value = "'0.07771"
$ws.Range("E5").Value = "'5.40%"
$ws.Range("D6").Value = "'2.357"
$ws.Range("E6").Value = "'4.47%"
$ws.Range("D7").Value = "'8.021"
$ws.Range("E7").Value = "'4.16%"
$ws.Range("D8").Value = "'3.936"
$ws.Range("E8").Value = "'5.41%"
$ws.Range("D9").Value = "'0.9263"
$ws.Range("D10").Value = "'0.1007"
$ws.Range("E10").Value = "'15.39%"
$ws.Range("E11").Value = "'6.51%"
$ws.Range("D12").Value = "'0.08513"
$ws.Range("E12").Value = "'4.38%"
$ws.Range("D13").Value = "'0.03316"
$ws.Range("E13").Value = "'6.50%"
$ws.Range("D14").Value = "'0.09893"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'-0.02%"
$ws.Range("D16").Value = "'0.005761"
$ws.Range("E16").Value = "'0.79%"
$ws.Range("D17").Value = "'3.470"
$ws.Range("E17").Value = "'-0.59%"
$ws.Range("E18").Value = "'3.90%"
$ws.Range("E19").Value = "'1.23%"
$ws.Range("D20").Value = "'0.1310"
$ws.Range("E20").Value = "'1.56%"
$ws.Range("D21").Value = "'4.291"
$ws.Range("E21").Value = "'12.12%"
$ws.Range("E22").Value = "'12.35%"
$ws.Range("D23").Value = "'0.04564"
$ws.Range("E23").Value = "'0.25%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.36%"
$ws.Range("D25").Value = "'0.004468"
$ws.Range("E25").Value = "'7.91%"
$ws.Range("E26").Value = "'-3.99%"
$ws.Range("E27").Value = "'8.82%"
$ws.Range("D39").Value = "'0.01783"
$ws.Range("E39").Value = "'12.73%"
$ws.Range("D40").Value = "'0.04749"
$ws.Range("E40").Value = "'6.38%"
$ws.Range("E41").Value = "'5.69%"
$ws.Range("D42").Value = "'0.1414"
$ws.Range("E42").Value = "'6.77%"
$ws.Range("D43").Value = "'0.007077"
$ws.Range("E43").Value = "'-26.05%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-6.59%"
$ws.Range("D45").Value = "'0.009523"
$ws.Range("E45").Value = "'12.92%"
$ws.Range("D46").Value = "'0.00006116"
$ws.Range("E46").Value = "'-0.33%"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("D48").Value = "'2.731"
$ws.Range("E48").Value = "'24.72%"
$ws.Range("D49").Value = "'0.001999"
$ws.Range("E49").Value = "'-0.17%"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.17%"
